{"js": "// Update the two-digit \u00f7 one-digit practice table: replace each\n// exercise's text with its new value, matching the commit's diff\n// (25 `<w:t>` runs inside the table change; the date heading is untouched).\n//\n// The 5x5 table's filled rows are 0, 4, 8, 12, 16 (the others are spacer\n// rows); each filled row has 5 cells. Addressing by (row, col) \u2014 rather than\n// by searching for the old text \u2014 avoids collisions, since some new values\n// equal other cells' old values (e.g. \"22\u00f74=\" becomes \"41\u00f74=\", while the\n// existing \"41\u00f74=\" cell separately becomes \"37\u00f77=\").\nconst rowNewValues = {\n  0: [\"67\u00f72=\", \"42\u00f79=\", \"45\u00f72=\", \"41\u00f74=\", \"21\u00f73=\"],\n  4: [\"44\u00f73=\", \"95\u00f76=\", \"33\u00f74=\", \"86\u00f73=\", \"29\u00f73=\"],\n  8: [\"24\u00f75=\", \"22\u00f76=\", \"82\u00f73=\", \"71\u00f75=\", \"94\u00f78=\"],\n  12: [\"15\u00f76=\", \"96\u00f72=\", \"10\u00f78=\", \"35\u00f77=\", \"13\u00f76=\"],\n  16: [\"36\u00f74=\", \"37\u00f77=\", \"41\u00f77=\", \"28\u00f79=\", \"66\u00f76=\"],\n};\n\nconst table = context.document.body.tables.getFirst();\n\nfor (const rowIndex of Object.keys(rowNewValues)) {\n  const r = Number(rowIndex);\n  const newValues = rowNewValues[r];\n  for (let c = 0; c < newValues.length; c++) {\n    const cell = table.getCell(r, c);\n    const range = cell.body.getRange();\n    range.insertText(newValues[c], \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the two-digit \u00f7 one-digit practice table: replace each\n# exercise's text with its new value, matching the commit's diff\n# (25 runs inside the table change; the date heading is untouched).\n#\n# The 5x5 table's filled rows are Word's 1-based rows 1, 5, 9, 13, 17 (the\n# rows in between are blank spacer rows); each filled row has 5 cells.\n# Addressing by (row, col) \u2014 rather than searching for the old text \u2014 avoids\n# collisions, since some new values equal other cells' old values (e.g.\n# \"22\u00f74=\" becomes \"41\u00f74=\", while the existing \"41\u00f74=\" cell separately\n# becomes \"37\u00f77=\").\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$rowNewValues = @{\n    1  = @(\"67\u00f72=\", \"42\u00f79=\", \"45\u00f72=\", \"41\u00f74=\", \"21\u00f73=\")\n    5  = @(\"44\u00f73=\", \"95\u00f76=\", \"33\u00f74=\", \"86\u00f73=\", \"29\u00f73=\")\n    9  = @(\"24\u00f75=\", \"22\u00f76=\", \"82\u00f73=\", \"71\u00f75=\", \"94\u00f78=\")\n    13 = @(\"15\u00f76=\", \"96\u00f72=\", \"10\u00f78=\", \"35\u00f77=\", \"13\u00f76=\")\n    17 = @(\"36\u00f74=\", \"37\u00f77=\", \"41\u00f77=\", \"28\u00f79=\", \"66\u00f76=\")\n}\n\nforeach ($row in $rowNewValues.Keys) {\n    $values = $rowNewValues[$row]\n    for ($col = 1; $col -le $values.Count; $col++) {\n        $cell = $t.Cell($row, $col)\n        $cell.Range.Text = $values[$col - 1]\n    }\n}\n"}
